$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.010.77"
$ws.Range("E2").Value = "  -0.72%  "

# Row 3
$ws.Range("D3").Value = "1.648.05"
$ws.Range("E3").Value = "  -0.53%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.52"
$ws.Range("E5").Value = "  -0.64%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5216"
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("E7").Value = "  -0.27%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2614"
$ws.Range("E8").Value = "  -2.14%  "

# Row 9
$ws.Range("E9").Value = "  -0.84%  "

# Row 10
$ws.Range("E10").Value = "  -3.05%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07738"
$ws.Range("E11").Value = "  -0.11%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.691.60"
$ws.Range("E12").Value = "  +2.02%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.464"
$ws.Range("E13").Value = "  +0.57%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5437"
$ws.Range("E14").Value = "  -0.83%  "

# Row 15
$ws.Range("D15").Value = "0.0₅8075"
$ws.Range("E15").Value = "  -2.07%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.64"
$ws.Range("E16").Value = "  -0.57%  "

# Row 17
$ws.Range("D17").Value = "26.025.83"

# Row 18
$ws.Range("E18").Value = "  -0.30%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.557"
$ws.Range("E19").Value = "  -2.68%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.47"
$ws.Range("E20").Value = "  -0.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"

# Row 22
$ws.Range("E22").Value = "  -2.30%  "

# Row 23
$ws.Range("E23").Value = "  -0.39%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "138.83"
$ws.Range("E24").Value = "  +0.88%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1231"
$ws.Range("E25").Value = "  -0.50%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.240"
$ws.Range("E26").Value = "  -0.40%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.15"
$ws.Range("E27").Value = "  +0.17%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.428"
$ws.Range("E28").Value = "  +0.90%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05910"
$ws.Range("E29").Value = "  -2.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.273"
$ws.Range("E30").Value = "  -1.10%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.489"
$ws.Range("E31").Value = "  -1.84%  "

# Row 32
$ws.Range("E32").Value = "  -3.74%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.512"
$ws.Range("E33").Value = "  -8.41%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.415"
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9417"
$ws.Range("E35").Value = "  -4.18%  "

# Row 36
$ws.Range("E36").Value = "  -1.10%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5676"
$ws.Range("E37").Value = "  -5.00%  "

# Row 38
$ws.Range("E38").Value = "  +0.55%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.840"
$ws.Range("E39").Value = "  -2.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8462"
$ws.Range("E40").Value = "  -2.44%  "

# Row 41
$ws.Range("E41").Value = "  -0.23%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.47"
$ws.Range("E42").Value = "  +0.63%  "

# Row 43
$ws.Range("D43").Value = "1.000.13"
$ws.Range("E43").Value = "  -4.22%  "

# Row 44
$ws.Range("D44").Value = "1.790.31"
$ws.Range("E44").Value = "  -0.45%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "56.57"
$ws.Range("E45").Value = "  -1.19%  "

# Row 46
$ws.Range("E46").Value = "  -2.28%  "

# Row 47
$ws.Range("E47").Value = "  -0.45%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4294"
$ws.Range("E48").Value = "  +1.49%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.854"
$ws.Range("E49").Value = "  -3.28%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05148"
$ws.Range("E50").Value = "  -0.67%  "

# Row 51
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.469"
$ws.Range("E51").Value = "  -0.67%  "
